$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748891365835/F7KsVXPM8Y-camera-video-segment-driver-1748891368335.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNBPAN2TQ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQC9gUtLE7ybiFqW4RI39PFfxsVhgD5PUwFVkjHbHwsf3wIhAPYz0%2BptvjlUO%2FTtzOcnuhraFpn5fRsQ91eR1%2FbnZo1hKt0DCB4QBBoMNzgxMjA0OTQyMjQ0Igxa11vWVDboXM1zdKsqugOyozAb%2BfTlAoB%2F5662phY68qGM8SQQwhc0u7dJztmnJ0Q0BmuWlMvzcNReS%2FIrFOdpdDJn2T0LZNYuCXxsTis7Pm%2FGXM0WmJdBbZxoX9AqTXl3tOwiDdSbRSYZDkS15BU5BJMNN0Du1xXVafKNDvVyP90%2Ba3YRMcb2Z0OMYmNPFOHfpeSSdZ7gEI4o7%2BUJWW%2FEtfW4K5IAch3fCkZOqiHVHjfXenwMeYbWbQPm2PhyZrZk7b02b0rog%2Bs2iQ7DUdIxKaf4FFm13xDUI%2Fc5QBsgOmFg61QIaQ%2Fln85BEyPt6224ZsDwS%2BzaXlDYPIkJCfcyaHbU6Q0IJt2RIi9bMZFKOVjHrrzzKFcEgff4kKG0IkdMZZU%2BAAxhXx%2Bs45BqytrsZk8i8lwSKOc2c6ClNT7daXHyPw7E2qTx4dXJHa7khXsavhKZHQN%2BAZwgbFkpHG5EzcBQqMU0Iuv8jAnzW4x2q0ppxUicnzRSD%2BocTcKgWjkDbXd7Std9vRP2Y1QJUfCTgzbNb%2BecBh6bL6fmnDrbrmOI85g%2BxBsOqkLgfeiEicd6c4fwqfE%2BgRgSvJ0WYyAKBVGQTTn8dBUTMOa%2B%2FcEGOqQBXMezqQNmZAQYY9mxW10dMNBXUBobr%2FJJVBnlGkdvcEEbn%2BnF7gBlyKhBsS44rL%2FblAaqsGZisjYQRhmVb4j6b0acDil3dINcHa6xxKaodQmk%2B6OXVf1fJZHX82wHKHjzxdyvhDy6MMETK0tWnXpsdAgAc40YB61j5pWGZTps6zeWq2I0gVP9qtN9fleP8pDmFDkS0EtnNF12VTDP1XYeKT%2FV0x8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=996f6cf5eb44e4953b200854919ef55e1af2838cafcf4f68b0a4e530df8cea34'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748891365835/MSujKlorgb-camera-video-segment-1748891368335.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNBPAN2TQ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQC9gUtLE7ybiFqW4RI39PFfxsVhgD5PUwFVkjHbHwsf3wIhAPYz0%2BptvjlUO%2FTtzOcnuhraFpn5fRsQ91eR1%2FbnZo1hKt0DCB4QBBoMNzgxMjA0OTQyMjQ0Igxa11vWVDboXM1zdKsqugOyozAb%2BfTlAoB%2F5662phY68qGM8SQQwhc0u7dJztmnJ0Q0BmuWlMvzcNReS%2FIrFOdpdDJn2T0LZNYuCXxsTis7Pm%2FGXM0WmJdBbZxoX9AqTXl3tOwiDdSbRSYZDkS15BU5BJMNN0Du1xXVafKNDvVyP90%2Ba3YRMcb2Z0OMYmNPFOHfpeSSdZ7gEI4o7%2BUJWW%2FEtfW4K5IAch3fCkZOqiHVHjfXenwMeYbWbQPm2PhyZrZk7b02b0rog%2Bs2iQ7DUdIxKaf4FFm13xDUI%2Fc5QBsgOmFg61QIaQ%2Fln85BEyPt6224ZsDwS%2BzaXlDYPIkJCfcyaHbU6Q0IJt2RIi9bMZFKOVjHrrzzKFcEgff4kKG0IkdMZZU%2BAAxhXx%2Bs45BqytrsZk8i8lwSKOc2c6ClNT7daXHyPw7E2qTx4dXJHa7khXsavhKZHQN%2BAZwgbFkpHG5EzcBQqMU0Iuv8jAnzW4x2q0ppxUicnzRSD%2BocTcKgWjkDbXd7Std9vRP2Y1QJUfCTgzbNb%2BecBh6bL6fmnDrbrmOI85g%2BxBsOqkLgfeiEicd6c4fwqfE%2BgRgSvJ0WYyAKBVGQTTn8dBUTMOa%2B%2FcEGOqQBXMezqQNmZAQYY9mxW10dMNBXUBobr%2FJJVBnlGkdvcEEbn%2BnF7gBlyKhBsS44rL%2FblAaqsGZisjYQRhmVb4j6b0acDil3dINcHa6xxKaodQmk%2B6OXVf1fJZHX82wHKHjzxdyvhDy6MMETK0tWnXpsdAgAc40YB61j5pWGZTps6zeWq2I0gVP9qtN9fleP8pDmFDkS0EtnNF12VTDP1XYeKT%2FV0x8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=2ad6194e159a27c9e975d5a045199de8a5d392c129c98dd310c715e41b9e3bc6'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748885144009/VaUCDIqpi8-camera-video-segment-driver-1748885146509.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNBPAN2TQ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQC9gUtLE7ybiFqW4RI39PFfxsVhgD5PUwFVkjHbHwsf3wIhAPYz0%2BptvjlUO%2FTtzOcnuhraFpn5fRsQ91eR1%2FbnZo1hKt0DCB4QBBoMNzgxMjA0OTQyMjQ0Igxa11vWVDboXM1zdKsqugOyozAb%2BfTlAoB%2F5662phY68qGM8SQQwhc0u7dJztmnJ0Q0BmuWlMvzcNReS%2FIrFOdpdDJn2T0LZNYuCXxsTis7Pm%2FGXM0WmJdBbZxoX9AqTXl3tOwiDdSbRSYZDkS15BU5BJMNN0Du1xXVafKNDvVyP90%2Ba3YRMcb2Z0OMYmNPFOHfpeSSdZ7gEI4o7%2BUJWW%2FEtfW4K5IAch3fCkZOqiHVHjfXenwMeYbWbQPm2PhyZrZk7b02b0rog%2Bs2iQ7DUdIxKaf4FFm13xDUI%2Fc5QBsgOmFg61QIaQ%2Fln85BEyPt6224ZsDwS%2BzaXlDYPIkJCfcyaHbU6Q0IJt2RIi9bMZFKOVjHrrzzKFcEgff4kKG0IkdMZZU%2BAAxhXx%2Bs45BqytrsZk8i8lwSKOc2c6ClNT7daXHyPw7E2qTx4dXJHa7khXsavhKZHQN%2BAZwgbFkpHG5EzcBQqMU0Iuv8jAnzW4x2q0ppxUicnzRSD%2BocTcKgWjkDbXd7Std9vRP2Y1QJUfCTgzbNb%2BecBh6bL6fmnDrbrmOI85g%2BxBsOqkLgfeiEicd6c4fwqfE%2BgRgSvJ0WYyAKBVGQTTn8dBUTMOa%2B%2FcEGOqQBXMezqQNmZAQYY9mxW10dMNBXUBobr%2FJJVBnlGkdvcEEbn%2BnF7gBlyKhBsS44rL%2FblAaqsGZisjYQRhmVb4j6b0acDil3dINcHa6xxKaodQmk%2B6OXVf1fJZHX82wHKHjzxdyvhDy6MMETK0tWnXpsdAgAc40YB61j5pWGZTps6zeWq2I0gVP9qtN9fleP8pDmFDkS0EtnNF12VTDP1XYeKT%2FV0x8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=baf448472f1e8f898d50a06f351ba4294e800324d1f3ba78c05ae0745614553e'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748885144009/8zw63V5M8l-camera-video-segment-1748885146509.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNBPAN2TQ%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQC9gUtLE7ybiFqW4RI39PFfxsVhgD5PUwFVkjHbHwsf3wIhAPYz0%2BptvjlUO%2FTtzOcnuhraFpn5fRsQ91eR1%2FbnZo1hKt0DCB4QBBoMNzgxMjA0OTQyMjQ0Igxa11vWVDboXM1zdKsqugOyozAb%2BfTlAoB%2F5662phY68qGM8SQQwhc0u7dJztmnJ0Q0BmuWlMvzcNReS%2FIrFOdpdDJn2T0LZNYuCXxsTis7Pm%2FGXM0WmJdBbZxoX9AqTXl3tOwiDdSbRSYZDkS15BU5BJMNN0Du1xXVafKNDvVyP90%2Ba3YRMcb2Z0OMYmNPFOHfpeSSdZ7gEI4o7%2BUJWW%2FEtfW4K5IAch3fCkZOqiHVHjfXenwMeYbWbQPm2PhyZrZk7b02b0rog%2Bs2iQ7DUdIxKaf4FFm13xDUI%2Fc5QBsgOmFg61QIaQ%2Fln85BEyPt6224ZsDwS%2BzaXlDYPIkJCfcyaHbU6Q0IJt2RIi9bMZFKOVjHrrzzKFcEgff4kKG0IkdMZZU%2BAAxhXx%2Bs45BqytrsZk8i8lwSKOc2c6ClNT7daXHyPw7E2qTx4dXJHa7khXsavhKZHQN%2BAZwgbFkpHG5EzcBQqMU0Iuv8jAnzW4x2q0ppxUicnzRSD%2BocTcKgWjkDbXd7Std9vRP2Y1QJUfCTgzbNb%2BecBh6bL6fmnDrbrmOI85g%2BxBsOqkLgfeiEicd6c4fwqfE%2BgRgSvJ0WYyAKBVGQTTn8dBUTMOa%2B%2FcEGOqQBXMezqQNmZAQYY9mxW10dMNBXUBobr%2FJJVBnlGkdvcEEbn%2BnF7gBlyKhBsS44rL%2FblAaqsGZisjYQRhmVb4j6b0acDil3dINcHa6xxKaodQmk%2B6OXVf1fJZHX82wHKHjzxdyvhDy6MMETK0tWnXpsdAgAc40YB61j5pWGZTps6zeWq2I0gVP9qtN9fleP8pDmFDkS0EtnNF12VTDP1XYeKT%2FV0x8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=b5300ab2466b8aff80442fec2b6b3e417463e2da9e903491eb5363126c16e9c1'
